$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 336.6316
$ws.Range("I6").Value = 149.75
$ws.Range("J6").Value = 1333.3334
$ws.Range("K6").Value = 449.25
$ws.Range("L6").Value = 4000.0002
$ws.Range("M6").Value = -337.25
$ws.Range("N6").Value = -4224.0002
$ws.Range("H98").Value = 3465.1714
$ws.Range("I98").Value = 1534.8214
$ws.Range("J98").Value = 11186.571
$ws.Range("K98").Value = 1534.8214
$ws.Range("L98").Value = 11186.571
$ws.Range("M98").Value = -36.82140000000004
$ws.Range("N98").Value = -14182.571
$ws.Range("H113").Value = 4358.909
$ws.Range("I113").Value = 4295
$ws.Range("J113").Value = 4529.3335
$ws.Range("K113").Value = 4295
$ws.Range("L113").Value = 4529.3335
$ws.Range("M113").Value = -1041
$ws.Range("N113").Value = -11037.3335
$ws.Range("H116").Value = 1638.25
$ws.Range("I116").Value = 1377.5
$ws.Range("J116").Value = 1725.1666
$ws.Range("K116").Value = 1377.5
$ws.Range("L116").Value = 1725.1666
$ws.Range("M116").Value = 2064.5
$ws.Range("N116").Value = -8609.1666
$ws.Range("H122").Value = 3465.1714
$ws.Range("I122").Value = 1534.8214
$ws.Range("J122").Value = 11186.571
$ws.Range("K122").Value = 4604.4642
$ws.Range("L122").Value = 33559.713
$ws.Range("M122").Value = -2154.4642
$ws.Range("N122").Value = -38459.713

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9318.870000000001
$ws.Range("I32").Value = 5920.5513
$ws.Range("J32").Value = 21367.455
$ws.Range("K32").Value = 5920.5513
$ws.Range("L32").Value = 21367.455
$ws.Range("M32").Value = -5633.5513
$ws.Range("N32").Value = -21941.455
$ws.Range("H101").Value = 29500
$ws.Range("J101").Value = 29500
$ws.Range("L101").Value = 29500
$ws.Range("N101").Value = -35990
$ws.Range("H110").Value = 1360.2307
$ws.Range("I110").Value = 1504.8889
$ws.Range("J110").Value = 1034.75
$ws.Range("K110").Value = 1504.8889
$ws.Range("L110").Value = 1034.75
$ws.Range("M110").Value = 540.1111000000001
$ws.Range("N110").Value = -5124.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1362.2
$ws.Range("I107").Value = 1163.3
$ws.Range("J107").Value = 1760
$ws.Range("K107").Value = 1163.3
$ws.Range("L107").Value = 1760
$ws.Range("M107").Value = 756.7
$ws.Range("N107").Value = -5600
$ws.Range("H133").Value = 39483.168
$ws.Range("J133").Value = 39483.168
$ws.Range("L133").Value = 39483.168
$ws.Range("N133").Value = -49603.168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2214.9434
$ws.Range("I31").Value = 1284.1666
$ws.Range("J31").Value = 2985.2415
$ws.Range("K31").Value = 1284.1666
$ws.Range("L31").Value = 2985.2415
$ws.Range("M31").Value = -989.1666
$ws.Range("N31").Value = -3575.2415
$ws.Range("H34").Value = 2214.9434
$ws.Range("I34").Value = 1284.1666
$ws.Range("J34").Value = 2985.2415
$ws.Range("K34").Value = 1284.1666
$ws.Range("L34").Value = 2985.2415
$ws.Range("M34").Value = -1082.1666
$ws.Range("N34").Value = -3389.2415
$ws.Range("H133").Value = 31665.2
$ws.Range("J133").Value = 31665.2
$ws.Range("L133").Value = 31665.2
$ws.Range("N133").Value = -36725.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2288
$ws.Range("I81").Value = 133.33333
$ws.Range("J81").Value = 4442.6665
$ws.Range("K81").Value = 399.99999
$ws.Range("L81").Value = 13327.9995
$ws.Range("M81").Value = 723.00001
$ws.Range("N81").Value = -15573.9995
$ws.Range("H84").Value = 2288
$ws.Range("I84").Value = 133.33333
$ws.Range("J84").Value = 4442.6665
$ws.Range("K84").Value = 1199.99997
$ws.Range("L84").Value = 39983.9985
$ws.Range("M84").Value = 4416.00003
$ws.Range("N84").Value = -51215.9985
$ws.Range("H113").Value = 842329.5
$ws.Range("I113").Value = 2331546.8
$ws.Range("J113").Value = 598
$ws.Range("K113").Value = 6994640.399999999
$ws.Range("L113").Value = 1794
$ws.Range("M113").Value = -6992470.399999999
$ws.Range("N113").Value = -6134
$ws.Range("H122").Value = 738.69446
$ws.Range("I122").Value = 488.92856
$ws.Range("J122").Value = 1612.875
$ws.Range("K122").Value = 4400.35704
$ws.Range("L122").Value = 14515.875
$ws.Range("M122").Value = -1950.35704
$ws.Range("N122").Value = -19415.875
$ws.Range("H129").Value = 2061.5476
$ws.Range("I129").Value = 1677.625
$ws.Range("J129").Value = 2297.8076
$ws.Range("K129").Value = 5032.875
$ws.Range("L129").Value = 6893.4228
$ws.Range("M129").Value = -32.875
$ws.Range("N129").Value = -16893.4228
$ws.Range("H131").Value = 877.04
$ws.Range("J131").Value = 968.3022999999999
$ws.Range("L131").Value = 2904.9069
$ws.Range("N131").Value = -12984.9069
$ws.Range("H138").Value = 1182.6818
$ws.Range("I138").Value = 1127.3684
$ws.Range("J138").Value = 1533
$ws.Range("K138").Value = 3382.1052
$ws.Range("L138").Value = 4599
$ws.Range("M138").Value = 1757.8948
$ws.Range("N138").Value = -14879

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1419.4445
$ws.Range("I122").Value = 1435
$ws.Range("K122").Value = 4305
$ws.Range("M122").Value = -1855
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 3644.7302
$ws.Range("I132").Value = 2074.5557
$ws.Range("J132").Value = 5738.2964
$ws.Range("K132").Value = 6223.6671
$ws.Range("L132").Value = 17214.8892
$ws.Range("M132").Value = -3693.6671
$ws.Range("N132").Value = -22274.8892
$ws.Range("H133").Value = 37784
$ws.Range("J133").Value = 37784
$ws.Range("L133").Value = 37784
$ws.Range("N133").Value = -47904

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1356.3334
$ws.Range("I7").Value = 1341.2222
$ws.Range("K7").Value = 1341.2222
$ws.Range("M7").Value = -1229.2222
$ws.Range("H40").Value = 2863.7646
$ws.Range("I40").Value = 2963.1428
$ws.Range("K40").Value = 2963.1428
$ws.Range("M40").Value = -2827.1428
$ws.Range("H68").Value = 19695.818
$ws.Range("I68").Value = 34450.668
$ws.Range("J68").Value = 1990
$ws.Range("K68").Value = 34450.668
$ws.Range("L68").Value = 1990
$ws.Range("M68").Value = -33701.668
$ws.Range("N68").Value = -3488
$ws.Range("H71").Value = 19695.818
$ws.Range("I71").Value = 34450.668
$ws.Range("J71").Value = 1990
$ws.Range("K71").Value = 172253.34
$ws.Range("L71").Value = 9950
$ws.Range("M71").Value = -168509.34
$ws.Range("N71").Value = -17438
$ws.Range("H126").Value = 1356.3334
$ws.Range("I126").Value = 1341.2222
$ws.Range("K126").Value = 4023.6666
$ws.Range("M126").Value = -1553.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 19500
$ws.Range("J39").Value = 19500
$ws.Range("L39").Value = 19500
$ws.Range("N39").Value = -20326
$ws.Range("H43").Value = 10000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 10000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -10298
$ws.Range("H100").Value = 1738.1428
$ws.Range("I100").Value = 908
$ws.Range("J100").Value = 2568.2856
$ws.Range("K100").Value = 1816
$ws.Range("L100").Value = 5136.5712
$ws.Range("M100").Value = -1275
$ws.Range("N100").Value = -6218.5712
$ws.Range("H122").Value = 108935
$ws.Range("I122").Value = 12164.375
$ws.Range("J122").Value = 302476.25
$ws.Range("K122").Value = 36493.125
$ws.Range("L122").Value = 907428.75
$ws.Range("M122").Value = -34043.125
$ws.Range("N122").Value = -912328.75
